$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) for rows 2-5 from serial 45208 (2023-10-09) to 45212 (2023-10-13)
$ws.Range("C2").Value2 = 45212
$ws.Range("C3").Value2 = 45212
$ws.Range("C4").Value2 = 45212
$ws.Range("C5").Value2 = 45212
